# StudentLoans.xlsx update: refresh data series from 1989-2013 to 1989-2016,
# rename "African American" -> "Black", convert Source/Note lines to rich
# text with bold lead-ins, add a new data row for 2016 to both tables, and
# touch up a handful of layout / page-setup details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------
# 1) Make room for the new "2016" rows.
#    Row 12 is currently a blank spacer row between the two tables;
#    insert a fresh row above it for the new debt-table 2016 data.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Insert()

# Row 24 (post-shift) is the blank spacer row below the percent table;
# insert a fresh row above it for the new percent-table 2016 data.
$ws.Rows.Item(25).Insert()

# ---------------------------------------------------------------------
# 2) Titles: "1989-2013" -> "1989-2016"
# ---------------------------------------------------------------------
$ws.Range("A1").Value2 = "Average Family Student Loan Debt for Those Age 25-55, 1989-2016"
$ws.Range("A14").Value2 = "Share of Families With Student Loan Debt for Those Ages 25-55, 1989-2016"

# ---------------------------------------------------------------------
# 3) Header row relabel: "African American" -> "Black" (both tables)
# ---------------------------------------------------------------------
$ws.Range("C2").Value2 = "Black"
$ws.Range("C15").Value2 = "Black"

# ---------------------------------------------------------------------
# 4) Debt table values (dollars) - refreshed figures, plus new 2016 row
# ---------------------------------------------------------------------
$debt = @{
    3  = @(1989, 1100.4069999999999, 1160.568, 897.58259999999996)
    4  = @(1992, 1321.3030000000001, 927.48239999999998, 793.06100000000004)
    5  = @(1995, 1885.5830000000001, 1922.886, 1396.864)
    6  = @(1998, 3115.2080000000001, 1325.2170000000001, 1436.579)
    7  = @(2001, 2979.72, 2223.7190000000001, 1729.8610000000001)
    8  = @(2004, 4035.6750000000002, 3987.1770000000001, 1651.79)
    9  = @(2007, 5263.567, 6111.0280000000002, 3005.4769999999999)
    10 = @(2010, 8041.9889999999996, 9510.1, 3089.1660000000002)
    11 = @(2013, 8363.6049999999996, 10302.66, 3177.41)
    12 = @(2016, 11108.41, 14224.77, 7493.9989999999998)
}
foreach ($r in $debt.Keys) {
    $vals = $debt[$r]
    $ws.Range("A$r").Value2 = $vals[0]
    $ws.Range("B$r").Value2 = $vals[1]
    $ws.Range("C$r").Value2 = $vals[2]
    $ws.Range("D$r").Value2 = $vals[3]
}

# Number format + alignment for the debt-table values: whole dollars, right-aligned
$ws.Range("B3:D12").HorizontalAlignment = -4152
$ws.Range("B3:D12").NumberFormat = "0"

# ---------------------------------------------------------------------
# 5) Blank spacer row (row 13) between the two tables - explicit cells
#    with left alignment (matches the row beneath the debt table).
# ---------------------------------------------------------------------
$ws.Range("B13:D13").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 6) Percent table values - refreshed figures, plus new 2016 row
# ---------------------------------------------------------------------
$pct = @{
    16 = @(1989, 0.10471229999999999, 0.1788198, 0.12725230000000001)
    17 = @(1992, 0.14333399999999999, 0.13827709999999999, 0.091457399999999994)
    18 = @(1995, 0.16394529999999999, 0.1696231, 0.1355392)
    19 = @(1998, 0.15646679999999999, 0.14063580000000001, 0.1022798)
    20 = @(2001, 0.1467881, 0.18491850000000001, 0.13408100000000001)
    21 = @(2004, 0.18294060000000001, 0.21600839999999999, 0.1035643)
    22 = @(2007, 0.1969195, 0.28346769999999999, 0.14205690000000001)
    23 = @(2010, 0.27978340000000002, 0.32128679999999998, 0.14434759999999999)
    24 = @(2013, 0.28455550000000002, 0.41227730000000001, 0.1570289)
    25 = @(2016, 0.33675110000000003, 0.41835879999999998, 0.21896889999999999)
}
foreach ($r in $pct.Keys) {
    $vals = $pct[$r]
    $ws.Range("A$r").Value2 = $vals[0]
    $ws.Range("B$r").Value2 = $vals[1]
    $ws.Range("C$r").Value2 = $vals[2]
    $ws.Range("D$r").Value2 = $vals[3]
}
$ws.Range("B25:D25").NumberFormat = "0.0%"

# ---------------------------------------------------------------------
# 7) Footer text: Source / Note lines become rich text with a bold
#    lead-in run; dates refreshed 2013 -> 2016.
# ---------------------------------------------------------------------
$ws.Range("A27").Value2 = "Source:" + $nbsp + "Urban Institute calculations from Survey of Consumer Finances 1989" + [char]0x2013 + "2016."
$ws.Range("A27").Characters(1, 7).Font.Bold = $true

$ws.Range("A28").Value2 = "Note:" + $nbsp + "2016 dollars. Age is defined as the age of the household head."
$ws.Range("A28").Characters(1, 6).Font.Bold = $true

# "For more, visit http://urbn.is/wealthcharts" text rides down with the row
# shift automatically, but the Hyperlinks collection does not retarget its
# backing range on its own - rebuild it explicitly at the new row (A29) and
# restore the usual "black lead-in + blue underlined URL" run formatting.
$ws.Hyperlinks.Delete()
$ws.Range("A29").Value2 = "For more, visit http://urbn.is/wealthcharts"
$null = $ws.Hyperlinks.Add($ws.Cells.Item(29, 1), "http://urbn.is/wealthcharts", [Type]::Missing, [Type]::Missing, "http://urbn.is/wealthcharts")

$leadIn = $ws.Range("A29").Characters(1, 16)
$leadIn.Font.Color = 0
$leadIn.Font.Underline = $false

$urlPart = $ws.Range("A29").Characters(17, 28)
$urlPart.Font.Color = 0xFF0000
$urlPart.Font.Underline = $true

# ---------------------------------------------------------------------
# 8) Worksheet-level layout touch-ups
# ---------------------------------------------------------------------
$null = $ws.Range("A21").Select()
$null = $ws.Range("B23").Select()

$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.Orientation = 2
$ws.PageSetup.Zoom = 88

Write-Output "done"
